# Replace the content of 8 phishing-message paragraphs in the questionnaire.
# Manual line breaks inside a Word paragraph are represented, when read
# through Range.Text, as Chr(11) (vertical-tab). Writing Chr(11) back into
# Range.Text re-creates the <w:br/> elements, so we build each new message
# by joining its text segments with [char]11.

$d = $word.ActiveDocument
$brk = [char]11

# --- Paragraph 8: Spotify message -> "Hello Mr Louise Frick" message ----
$p8 = @(
    "Hello Mr Louise Frick,",
    "Always wanted to go to the São Paulo International Film Festival (Mostra Internacional de Cinema)? This is your chance!",
    "The Ministry of Culture of Santana de Parnaíba is giving away a number of tickets for the festival that will take place this October and we would be more than happy to see you.",
    "If you are interested, please reply to this message with your full details including address, credit card and ID.",
    "The number of places is limited!"
) -join $brk
$d.Paragraphs(8).Range.Text = $p8

# --- Paragraph 11: Supreme Pets message -> "Dear Gary  Leal" fitness tracker message ---
$p11 = @(
    "Dear Gary  Leal",
    "We've noticed unusual activity on your fitness tracker account. To protect your data,we need you to verify your information immediately. ",
    "Click the link below to confirm your identity. ",
    "Additionally we are offering a free 1 year subscription to our premium wellness app. Don't miss out! Verify your account now!"
) -join $brk
$d.Paragraphs(11).Range.Text = $p11

# --- Paragraph 15: "Greetings Mr Liddle" -> "Dear Vincent Cummins" account-locked message ---
$p15 = @(
    " Dear Vincent Cummins,",
    "    We are writing to inform you that your account has been locked due to a suspected security breach. ",
    "    In order to protect your account, we require you to provide the following information:",
    "    ",
    "    Credit card number:",
    "    Expiration date:",
    "    CVV code:",
    "    4-digit PIN:",
    "    ",
    "    Once we have received this information, we will be able to unlock your account and allow you to access your account again.",
    "    ",
    "    Thank you for your understanding and cooperation.",
    "    ",
    "    Best regards,",
    "    The customer service team"
) -join $brk
$d.Paragraphs(15).Range.Text = $p15

# --- Paragraph 17: "Dear Mr. Lindell" Apple message -> "Dear Mr. Cummins" turtle-donation message ---
$p17 = @(
    "Dear Mr. Cummins,",
    "We are approaching you since we know your fondness for all types of animals. Unfortunately, there has been a great catastrophe on the Island of Samora, where an oil tanker spilled, and all the sea turtles are in danger. We would greatly appreciate a 10-dollar donation to help save the poor turtles. ",
    "To donate, please send us your credit card details.",
    "",
    "The turtels will all thanks you!"
) -join $brk
$d.Paragraphs(17).Range.Text = $p17

# --- Paragraph 22: "Dear Nicole" -> "Hello Helen!" gym membership message ---
$p22 = "Hello Helen! A payment for your gym membership has failed. Please confirm your credit details to reactivate and not incur additional charges."
$d.Paragraphs(22).Range.Text = $p22

# --- Paragraph 24: "Dear Ms. Morrow" -> Qantas Airways message ---
$p24 = @(
    "Subject: Important Update from Qantas Airways",
    " ",
    " Dear Helen,",
    " ",
    " I hope this message finds you well. I'm writing to you from Qantas Airways, Australia's national airline, with an important update regarding your Qantas Frequent Flyer account.",
    " ",
    " Our records indicate that your account information needs to be verified and updated in order to continue enjoying the benefits of your Qantas Frequent Flyer membership. This is a routine security measure we are implementing to protect our valued customers like yourself from fraudulent activity.",
    " ",
    " To complete the verification process, please reply to this email with the following information:",
    " ",
    " - Full name",
    " - Date of birth",
    " - Credit card number",
    " - Expiration date",
    " - Security code (CVV)",
    " ",
    " Once we have verified your account details, you will be able to continue using your Qantas Frequent Flyer account without any interruption. We appreciate your prompt attention to this matter.",
    " ",
    " Thank you for your continued loyalty to Qantas Airways. We look forward to serving you on your next journey.",
    " ",
    " Sincerely,",
    " Qantas Airways Customer Support"
) -join $brk
$d.Paragraphs(24).Range.Text = $p24

# --- Paragraph 29: Beauty-offer message -> "Hi Christina," bank-loan message ---
$p29 = @(
    " Hi Christina,",
    "    Thank you for your application for a loan at the bank. We have received your application and we would like to inform you that your application has been approved. Your loan has been approved for 3000 TL. You can use this loan for your needs.",
    "    Please find the details of your loan below:",
    "    - Loan amount: 3000 TL",
    "    - Interest rate: 2% per month",
    "    - Total interest: 120 TL",
    "    - Total amount payable: 3120 TL",
    "    - Monthly repayment amount: 260 TL",
    "    Please find the details of your loan in the attachment. Please read the loan agreement carefully and sign the agreement. Please return the signed agreement to us within 7 days.",
    "    Thank you for your interest in our bank.",
    "    Best regards,",
    "    Bank",
    "    www.bank.com"
)
$p29text = $p29[0] + $brk + $brk + $p29[1] + $brk + $brk + $p29[2] + $brk + $brk + `
    $p29[3] + $brk + $p29[4] + $brk + $p29[5] + $brk + $p29[6] + $brk + $p29[7] + $brk + $brk + `
    $p29[8] + $brk + $brk + $p29[9] + $brk + $brk + $p29[10] + $brk + $p29[11] + $brk + $p29[12]
$d.Paragraphs(29).Range.Text = $p29text

# --- Paragraph 31: LaMer skincare message -> "Hello Christina Craft," blackmail message ---
$p31 = @(
    "Hello Christina Craft, ",
    "",
    "We have records of your illicit and explicit activities online. Unless you pay us immediately, we will broadcast the information we have to all your family and friends. Send your credit card information now. "
) -join $brk
$d.Paragraphs(31).Range.Text = $p31

Write-Output "done"
